$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# D21 currently holds the text-typed "2000" (확정공모가) for 디비금융스팩12호.
# That row is being replaced, but the same text value is needed at D18 for
# 한국스팩14호 in the refreshed dataset. Copy it now (as a value-only paste)
# so it keeps its original text representation instead of being re-typed as a
# literal string, which Excel would auto-convert to a number.
$ws.Range("D21").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = "아이빔테크놀로지"
$ws.Range("B2").Value = "2024.07.15~07.19"
$ws.Range("C2").Value = "7,300~8,500"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 16308
$ws.Range("F2").Value = "삼성증권"

$ws.Range("A3").Value = "이베스트스팩6호"
$ws.Range("B3").Value = "2024.06.27~06.28"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 8000
$ws.Range("F3").Value = "엘에스증권"

$ws.Range("A4").Value = "엑셀세라퓨틱스"
$ws.Range("B4").Value = "2024.06.24~06.28"
$ws.Range("C4").Value = "6,200~7,700"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 10032
$ws.Range("F4").Value = "대신증권"

$ws.Range("A5").Value = "피앤에스미캐닉스"
$ws.Range("B5").Value = "2024.06.24~06.28"
$ws.Range("C5").Value = "14,000~17,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 18900
$ws.Range("F5").Value = "키움증권"

$ws.Range("A6").Value = "이엔셀"
$ws.Range("B6").Value = "2024.06.17~06.21"
$ws.Range("C6").Value = "13,600~15,300"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 21308
$ws.Range("F6").Value = "NH투자증권"

$ws.Range("A7").Value = "이노그리드"
$ws.Range("B7").Value = "2024.06.13~06.19"
$ws.Range("C7").Value = "29,000~35,000"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 17400
$ws.Range("F7").Value = "한국투자증권"

$ws.Range("A8").Value = "하스"
$ws.Range("B8").Value = "2024.06.13~06.19"
$ws.Range("C8").Value = "9,000~12,000"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = 16290
$ws.Range("F8").Value = "삼성증권"

$ws.Range("A9").Value = "에이치브이엠(구.한국진공야금)"
$ws.Range("B9").Value = "2024.06.11~06.17"
$ws.Range("C9").Value = "11,000~14,200"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = 26400
$ws.Range("F9").Value = "NH투자증권"

$ws.Range("A10").Value = "이노스페이스"
$ws.Range("B10").Value = "2024.06.11~06.17"
$ws.Range("C10").Value = "36,400~43,300"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = 48412
$ws.Range("F10").Value = "미래에셋증권,신한투자증권"

$ws.Range("A11").Value = "한국스팩15호"
$ws.Range("B11").Value = "2024.06.10~06.11"
$ws.Range("C11").Value = "2,000~2,000"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = 12500
$ws.Range("F11").Value = "한국투자증권"

$ws.Range("A12").Value = "하이젠알앤엠"
$ws.Range("B12").Value = "2024.06.07~06.13"
$ws.Range("C12").Value = "4,500~5,500"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = 15300
$ws.Range("F12").Value = "한국투자증권"

$ws.Range("A13").Value = "미래에셋비전스팩6호"
$ws.Range("B13").Value = "2024.06.05~06.07"
$ws.Range("C13").Value = "2,000~2,000"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = 12900
$ws.Range("F13").Value = "미래에셋증권"

$ws.Range("A14").Value = "KB스팩29호"
$ws.Range("B14").Value = "2024.06.04~06.05"
$ws.Range("C14").Value = "2,000~2,000"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = 12000
$ws.Range("F14").Value = "KB증권"

$ws.Range("A15").Value = "에이치엠씨아이비스팩7호"
$ws.Range("B15").Value = "2024.06.04~06.05"
$ws.Range("C15").Value = "2,000~2,000"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = 14000
$ws.Range("F15").Value = "현대차증권"

$ws.Range("A16").Value = "에스오에스랩"
$ws.Range("B16").Value = "2024.06.03~06.10"
$ws.Range("C16").Value = "7,500~9,000"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = 15000
$ws.Range("F16").Value = "한국투자증권"

$ws.Range("A17").Value = "미래에셋비전스팩5호"
$ws.Range("B17").Value = "2024.06.03~06.04"
$ws.Range("C17").Value = "2,000~2,000"
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = 9500
$ws.Range("F17").Value = "미래에셋증권"

$ws.Range("A18").Value = "한국스팩14호"
$ws.Range("B18").Value = "2024.06.03~06.04"
$ws.Range("C18").Value = "2,000~2,000"
# D18 already set above via PasteSpecial (preserves text "2000")
$ws.Range("E18").Value = 8000
$ws.Range("F18").Value = "한국투자증권"

$ws.Range("A19").Value = "시프트업(유가)"
$ws.Range("B19").Value = "2024.06.03~06.27"
$ws.Range("C19").Value = "47,000~60,000"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = 340750
$ws.Range("F19").Value = "한국투자증권,NH투자증권,신한투자증권"

$ws.Range("A20").Value = "신한글로벌액티브"
$ws.Range("B20").Value = "2024.06.03~06.05"
$ws.Range("C20").Value = "3,000~3,800"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = 70000
$ws.Range("F20").Value = "신한투자증권,한국투자증권"

$ws.Range("A21").Value = "한중엔시에스"
$ws.Range("B21").Value = "2024.05.30~06.05"
$ws.Range("C21").Value = "20,000~23,500"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = 32000
$ws.Range("F21").Value = "IBK투자증권"
